$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B10").Value = "Amatenango De La Frontera"
$ws.Range("B11").Value = "Bejucal De Ocampo"
$ws.Range("B15").Value = "Comitán De Domínguez"
$ws.Range("B27").Value = "Mazapa De Madero"
$ws.Range("B47").Value = "Hidalgo Del Parral"
$ws.Range("A49").Value = "Ciudad De México"
$ws.Range("B62").Value = "San Juan De Guadalupe"
$ws.Range("A64").Value = "Estado De México"
$ws.Range("B64").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B66").Value = "Atizapán De Zaragoza"
$ws.Range("B69").Value = "Ecatepec De Morelos"
$ws.Range("B70").Value = "Ixtapan De La Sal"
$ws.Range("B71").Value = "Ixtapan Del Oro"
$ws.Range("B73").Value = "Naucalpan De Juárez"
$ws.Range("D74").Value = 0.009538950715421305
$ws.Range("B75").Value = "San Felipe Del Progreso"
$ws.Range("B79").Value = "Tlalnepantla De Baz"
$ws.Range("B83").Value = "Apaseo El Alto"
$ws.Range("D84").Value = 0.009538950715421305
$ws.Range("B86").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B91").Value = "San Luis De La Paz"
$ws.Range("B93").Value = "Acapulco De Juárez"
$ws.Range("B94").Value = "Ajuchitlán Del Progreso"
$ws.Range("B96").Value = "Ayutla De Los Libres"
$ws.Range("B99").Value = "Coyuca De Benítez"
$ws.Range("B100").Value = "Coyuca De Catalán"
$ws.Range("B101").Value = "Cutzamala De Pinzón"
$ws.Range("B105").Value = "Tlapa De Comonfort"
$ws.Range("B110").Value = "Molango De Escamilla"
$ws.Range("B111").Value = "Pachuca De Soto"
$ws.Range("B113").Value = "Tula De Allende"
$ws.Range("B117").Value = "Autlán De Navarro"
$ws.Range("B120").Value = "Concepción De Buenos Aires"
$ws.Range("B121").Value = "Cuautitlán De García Barragán"
$ws.Range("B124").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B125").Value = "Ixtlahuacán Del Río"
$ws.Range("B127").Value = "Lagos De Moreno"
$ws.Range("B129").Value = "Ojuelos De Jalisco"
$ws.Range("B134").Value = "Teocuitatlán De Corona"
$ws.Range("B136").Value = "Tizapán El Alto"
$ws.Range("B137").Value = "Tlajomulco De Zúñiga"
$ws.Range("B140").Value = "Unión De Tula"
$ws.Range("B163").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B172").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B174").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B175").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B177").Value = "Oaxaca De Juárez"
$ws.Range("B178").Value = "Putla Villa De Guerrero"
$ws.Range("B183").Value = "San Dionisio Del Mar"
$ws.Range("B211").Value = "Izúcar De Matamoros"
$ws.Range("B225").Value = "Tepexi De Rodríguez"
$ws.Range("B227").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B244").Value = "Santa María Del Río"
$ws.Range("B246").Value = "Villa De Reyes"
$ws.Range("D246").Value = 0.009538950715421305
$ws.Range("D265").Value = 0.009538950715421305
$ws.Range("B266").Value = "Contla De Juan Cuamatzi"
$ws.Range("B273").Value = "Camarón De Tejeda"
$ws.Range("B277").Value = "Ignacio De La Llave"
$ws.Range("B278").Value = "Ixhuatlán De Madero"
$ws.Range("B283").Value = "Martínez De La Torre"
$ws.Range("B288").Value = "Paso De Ovejas"
$ws.Range("B290").Value = "Sayula De Alemán"
$ws.Range("B291").Value = "Soledad De Doblado"
$ws.Range("B302").Value = "Tlaltenango De Sánchez Román"

# Remove the trailing metadata rows (305-310) that are no longer part of the dataset
$ws.Range("A305:A310").EntireRow.Delete()
